$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 11:04"

# Update country data rows with refreshed case counts.
# Some rows also change country (reordered) because the sheet is kept
# sorted by "Casos totales" (column B) descending; values are written
# directly per-row/column to reproduce the exact final layout.

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 3695469
$ws.Range("C4").Value = 444
$ws.Range("D4").Value = 1680418
$ws.Range("E4").Value = 1873927
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 141124

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 1005871
$ws.Range("C6").Value = 234
$ws.Range("D6").Value = 636727
$ws.Range("E6").Value = 343522
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 13
$ws.Range("H6").Value = 25622

# Row 20: Banglades
$ws.Range("A20").Value = "Banglades"
$ws.Range("B20").Value = 199357
$ws.Range("C20").Value = 3034
$ws.Range("D20").Value = 108725
$ws.Range("E20").Value = 88085
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 51
$ws.Range("H20").Value = 2547

# Row 35: Filipinas
$ws.Range("A35").Value = "Filipinas"
$ws.Range("B35").Value = 63001
$ws.Range("C35").Value = 1841
$ws.Range("D35").Value = 21748
$ws.Range("E35").Value = 39593
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 17
$ws.Range("H35").Value = 1660

# Row 36: Oman
$ws.Range("A36").Value = "Oman"
$ws.Range("B36").Value = 62574
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 40090
$ws.Range("E36").Value = 22194
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 290

# Row 45: Singapur
$ws.Range("A45").Value = "Singapur"
$ws.Range("B45").Value = 47453
$ws.Range("C45").Value = 327
$ws.Range("D45").Value = 43256
$ws.Range("E45").Value = 4170
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 27

# Row 46: Israel
$ws.Range("A46").Value = "Israel"
$ws.Range("B46").Value = 46546
$ws.Range("C46").Value = 487
$ws.Range("D46").Value = 20523
$ws.Range("E46").Value = 25636
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = 387

# Row 47: Polonia
$ws.Range("A47").Value = "Polonia"
$ws.Range("B47").Value = 39407
$ws.Range("C47").Value = 353
$ws.Range("D47").Value = 29505
$ws.Range("E47").Value = 8290
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 7
$ws.Range("H47").Value = 1612

# Row 48: Afganistan
$ws.Range("A48").Value = "Afganistan"
$ws.Range("B48").Value = 35229
$ws.Range("C48").Value = 159
$ws.Range("D48").Value = 23151
$ws.Range("E48").Value = 10931
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 34
$ws.Range("H48").Value = 1147

# Row 49: Barein
$ws.Range("A49").Value = "Barein"
$ws.Range("B49").Value = 35084
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 30809
$ws.Range("E49").Value = 4154
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 121

# Row 80: Malasia
$ws.Range("A80").Value = "Malasia"
$ws.Range("B80").Value = 8755
$ws.Range("C80").Value = 18
$ws.Range("D80").Value = 8541
$ws.Range("E80").Value = 92
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 122

# Row 111: Sri Lanka
$ws.Range("A111").Value = "Sri Lanka"
$ws.Range("B111").Value = 2687
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 2012
$ws.Range("E111").Value = 664
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 11

# Row 118: Montenegro
$ws.Range("A118").Value = "Montenegro"
$ws.Range("B118").Value = 1965
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 357
$ws.Range("E118").Value = 1582
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 26

# Row 119: Eslovaquia
$ws.Range("A119").Value = "Eslovaquia"
$ws.Range("B119").Value = 1965
$ws.Range("C119").Value = 14
$ws.Range("D119").Value = 1523
$ws.Range("E119").Value = 414
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 28

# Row 120: Eslovenia
$ws.Range("A120").Value = "Eslovenia"
$ws.Range("B120").Value = 1916
$ws.Range("C120").Value = 19
$ws.Range("D120").Value = 1522
$ws.Range("E120").Value = 283
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 111

# Row 121: Islandia
$ws.Range("A121").Value = "Islandia"
$ws.Range("B121").Value = 1914
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 1892
$ws.Range("E121").Value = 12
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 10

# Row 122: Lituania
$ws.Range("A122").Value = "Lituania"
$ws.Range("B122").Value = 1908
$ws.Range("C122").Value = 6
$ws.Range("D122").Value = 1595
$ws.Range("E122").Value = 234
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 79

# Row 123: Guinea-Bisau
$ws.Range("A123").Value = "Guinea-Bisau"
$ws.Range("B123").Value = 1902
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 773
$ws.Range("E123").Value = 1103
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 26

# Row 124: Zambia
$ws.Range("A124").Value = "Zambia"
$ws.Range("B124").Value = 1895
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 1412
$ws.Range("E124").Value = 441
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 42

# Row 125: Cabo Verde
$ws.Range("A125").Value = "Cabo Verde"
$ws.Range("B125").Value = 1894
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 902
$ws.Range("E125").Value = 973
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 19

# Row 126: Hong Kong
$ws.Range("A126").Value = "Hong Kong"
$ws.Range("B126").Value = 1714
$ws.Range("C126").Value = 58
$ws.Range("D126").Value = 1264
$ws.Range("E126").Value = 440
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 10

# Row 127: Sierra Leona
$ws.Range("A127").Value = "Sierra Leona"
$ws.Range("B127").Value = 1678
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 1213
$ws.Range("E127").Value = 401
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 64

# Row 128: Libia
$ws.Range("A128").Value = "Libia"
$ws.Range("B128").Value = 1652
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 379
$ws.Range("E128").Value = 1227
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 46

# Row 129: Suazilandia
$ws.Range("A129").Value = "Suazilandia"
$ws.Range("B129").Value = 1552
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 736
$ws.Range("E129").Value = 795
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 21

# Row 130: Yemen
$ws.Range("A130").Value = "Yemen"
$ws.Range("B130").Value = 1552
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 695
$ws.Range("E130").Value = 419
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 438

# Row 131: Nueva Zelanda
$ws.Range("A131").Value = "Nueva Zelanda"
$ws.Range("B131").Value = 1549
$ws.Range("C131").Value = 1
$ws.Range("D131").Value = 1506
$ws.Range("E131").Value = 21
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 22

# Row 132: Ruanda
$ws.Range("A132").Value = "Ruanda"
$ws.Range("B132").Value = 1473
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 770
$ws.Range("E132").Value = 699
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 4

# Row 133: Benin
$ws.Range("A133").Value = "Benin"
$ws.Range("B133").Value = 1463
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 557
$ws.Range("E133").Value = 878
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 28

# Row 134: Mozambique
$ws.Range("A134").Value = "Mozambique"
$ws.Range("B134").Value = 1383
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 375
$ws.Range("E134").Value = 999
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 9

# Row 135: Zimbabue
$ws.Range("A135").Value = "Zimbabue"
$ws.Range("B135").Value = 1362
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 425
$ws.Range("E135").Value = 914
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 23

# Row 136: Tunez
$ws.Range("A136").Value = "Tunez"
$ws.Range("B136").Value = 1327
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 1093
$ws.Range("E136").Value = 184
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 50
